$d = $word.ActiveDocument

$wns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# Locate the "Liquor 白酒" paragraph (last real content paragraph, right
# before the trailing empty bookmark paragraph) and insert the three new
# vocabulary entries directly after it, one at a time, so each new
# paragraph lands just before the final empty paragraph.

function Insert-ParaAfterLast([string]$innerXml) {
    $count = $d.Paragraphs.Count
    $lastContentPara = $d.Paragraphs.Item($count - 1)
    $insertPt = $d.Range($lastContentPara.Range.End - 1, $lastContentPara.Range.End - 1)
    $xml = "<w:p $wns>$innerXml</w:p>"
    [void]$insertPt.InsertXML($xml)
}

# Paragraph 1: Stab 插
Insert-ParaAfterLast('<w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t xml:space="preserve">Stab </w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/><w:lang w:eastAsia="zh-HK"/></w:rPr><w:t>插</w:t></w:r>')

# Paragraph 2: Prophecy 預言
Insert-ParaAfterLast('<w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>P</w:t></w:r><w:r><w:t xml:space="preserve">rophecy </w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/><w:lang w:eastAsia="zh-HK"/></w:rPr><w:t>預言</w:t></w:r>')

# Paragraph 3: Prophet 預言家
Insert-ParaAfterLast('<w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t xml:space="preserve">Prophet </w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/><w:lang w:eastAsia="zh-HK"/></w:rPr><w:t>預言家</w:t></w:r>')

Write-Output "Inserted 3 new paragraphs. Total paragraphs now: $($d.Paragraphs.Count)"
